$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E3").Value = "UNIQUE"
$ws.Range("C3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
